$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml) - column F ("想去人数") value updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2746
$ws1.Range("F7").Value = 2388
$ws1.Range("F8").Value = 1854
$ws1.Range("F9").Value = 221
$ws1.Range("F11").Value = 2507
$ws1.Range("F12").Value = 564
$ws1.Range("F13").Value = 249
$ws1.Range("F16").Value = 131
$ws1.Range("F18").Value = 9348
$ws1.Range("F20").Value = 7262
$ws1.Range("F21").Value = 11827
$ws1.Range("F25").Value = 367
$ws1.Range("F27").Value = 2644
$ws1.Range("F29").Value = 203
$ws1.Range("F30").Value = 2599
$ws1.Range("F31").Value = 793
$ws1.Range("F33").Value = 4536
$ws1.Range("F34").Value = 974
$ws1.Range("F37").Value = 542

# Sheet "全部类型" (sheet4.xml) - column F ("想去人数") value updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 2746
$ws4.Range("F11").Value = 2388
$ws4.Range("F13").Value = 1854
$ws4.Range("F14").Value = 221
$ws4.Range("F15").Value = 2507
$ws4.Range("F17").Value = 564
$ws4.Range("F18").Value = 249
$ws4.Range("F21").Value = 131
$ws4.Range("F23").Value = 9348
$ws4.Range("F24").Value = 58
$ws4.Range("F25").Value = 7262
$ws4.Range("F26").Value = 11827
$ws4.Range("F30").Value = 367
$ws4.Range("F34").Value = 2644
$ws4.Range("F38").Value = 203
$ws4.Range("F40").Value = 4536
$ws4.Range("F45").Value = 542
